$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.607.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.26%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.369.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.95%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'557.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.28%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'176.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.98%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.37%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.358.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.12%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.14%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.630"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.98%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +2.70%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'54.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.69%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +0.39%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  -0.22%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.903.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.04%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'18.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.73%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -1.86%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.370.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.79%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'11.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.18%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'64.506.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.42%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  -0.35%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'459.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +12.94%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'4.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +9.72%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -1.78%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'85.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.17%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'13.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.23%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'10.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.50%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +1.82%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'8.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.11%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'30.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.60%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'6.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.99%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'11.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.11%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'584.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.26%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -0.49%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'58.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.59%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +0.16%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.140"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -8.10%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'35.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.89%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.12%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'" + "0.0" + [char]0x2083 + "0759"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.05%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.372"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.01%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'3.111.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.26%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -0.28%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.00%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.32%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +0.53%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -0.44%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -0.03%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -1.08%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'8.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.03%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'135.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.14%  "
$ws.Range("E51").Style = "Normal"
